# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> used by the Notes Master  (was "Office Theme")
#   ppt/theme/theme2.xml  -> used by the Slide Master / presentation ("Integral")
#
# The authored commit swaps the two parts' contents in place (file names /
# relationships are untouched, only the colour definitions move): the Slide
# Master's theme becomes the plain "Office Theme" palette that used to live
# in theme1.xml, while the Notes Master's theme becomes the "Integral"
# palette that used to live in theme2.xml.
#
# This COM host models a single shared Theme object (reached from
# SlideMaster / NotesMaster / HandoutMaster / Designs alike) that is backed
# by the Slide Master's theme part, so we recolor it here to the target
# "Office Theme" scheme.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# Index : scheme slot : target "Office Theme" RGB (packed as PowerPoint's
# BGR long, i.e. what ThemeColorScheme.Colors(i).RGB expects/returns)
$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
